$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I (shifts old I..L to J..M)
$ws.Columns.Item(9).Insert()

# --- Header row ---
$ws.Range("H1").Value = "Mode de paiement commandes"
$ws.Range("I1").Value = "Mode de paiement livraison"
$ws.Range("J1").Value = "Infos Pizzas"
$ws.Range("K1").Value = "Prix Livraison"
$ws.Range("L1").Value = "Prix Pizzas + Suppléments"
$ws.Range("M1").Value = "Total"

# --- Row 2 ---
$ws.Range("B2").Value = "Oui"
$ws.Range("C2").Value = "Commandée sur place"
$ws.Range("D2").Value = "Heure sur place: 17:40"
$ws.Range("E2").Value = "Commandée sur place"
$ws.Range("F2").Value = "Sur place"
$ws.Range("G2").Value = "Lahat Samb"
$ws.Range("I2").Value = "Chez le livreur"
$ws.Range("J2").Value = "Selucy Taille Petite, Margherita - Vosgienne Taille Grande Supplements: Fromages, Emmental"
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 9600
$ws.Range("M2").Value = 9600

# --- Row 3 ---
$ws.Range("B3").Value = "Non"
$ws.Range("C3").Value = "Ouest Foire Dakar"
$ws.Range("D3").Value = "19h53"
$ws.Range("E3").Value = "Cagil"
$ws.Range("F3").Value = "Livré"
$ws.Range("G3").Value = "George N'gock"
$ws.Range("H3").Value = "Chez Izoua"
$ws.Range("I3").Value = "Chez Izoua"
$ws.Range("J3").Value = "Selucy Taille Grande, Margherita - Vosgienne Taille Grande"
$ws.Range("K3").Value = "Commandée sur place"
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 11800

# --- Row 4 ---
$ws.Range("C4").Value = "Keur Gorgui"
$ws.Range("D4").Value = "22h55"
$ws.Range("E4").Value = "Bazoungoula"
$ws.Range("G4").Value = "Coulibaly Yelanto"
$ws.Range("H4").Value = "Chez Izoua"
$ws.Range("I4").Value = "Chez Izoua"
$ws.Range("J4").Value = "Vosgienne Taille Petite, Margherita - Selucy Taille Grande"
$ws.Range("K4").Value = "Commandée sur place"
$ws.Range("L4").Value = 7500
$ws.Range("M4").Value = 9500

# --- Row 5 (new) ---
$ws.Range("A5").Value = "27 décembre 2024"
$ws.Range("B5").Value = "Non"
$ws.Range("C5").Value = "Ouest Foire"
$ws.Range("D5").Value = "22h56"
$ws.Range("E5").Value = "Bazoungoula"
$ws.Range("F5").Value = "Annulé"
$ws.Range("G5").Value = "Yves"
$ws.Range("H5").Value = "Chez Izoua"
$ws.Range("I5").Value = "Chez Izoua"
$ws.Range("J5").Value = "Vosgienne - Selucy Taille Grande Supplements: Fromages, Emmental"
$ws.Range("K5").Value = "Commandée sur place"
$ws.Range("L5").Value = 7100
$ws.Range("M5").Value = 9100

# --- Row 6 (new) ---
$ws.Range("A6").Value = "27 décembre 2024"
$ws.Range("B6").Value = "Non"
$ws.Range("C6").Value = "Grand Dakar"
$ws.Range("D6").Value = "20h55"
$ws.Range("E6").Value = "Bazoungoula"
$ws.Range("F6").Value = "En attente"
$ws.Range("G6").Value = "Alphonse Desire"
$ws.Range("H6").Value = "Chez Izoua"
$ws.Range("I6").Value = "Chez Izoua"
$ws.Range("J6").Value = "Selucy Taille Grande, Vosgienne - Margherita Taille Grande"
$ws.Range("K6").Value = "Commandée sur place"
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = 11500
